# Update equipment status from "Disponible" (Available) to "Prestado" (Loaned)
# for the equipment rows that have been checked out: EQ-001, EQ-002, EQ-007,
# EQ-009, EQ-010 (rows 2, 3, 8, 10, 11 in the "Estado" column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2, 3, 8, 10, 11)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Prestado"
}
